# BUG: Fix read_excel w/parse_cols & empty dataset (#23661)
# Adds a third sheet ("Sheet3") containing only a header row (A..F) and
# makes it the active sheet/tab, matching the fixture used to reproduce
# gh-9208 (parsing an all-header, no-data sheet).

$wb = $excel.ActiveWorkbook

# Add a new worksheet at the end of the workbook and rename it to "Sheet3".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Header-only row: A, B, C, D, E, F
$ws3.Range("A1").Value = "A"
$ws3.Range("B1").Value = "B"
$ws3.Range("C1").Value = "C"
$ws3.Range("D1").Value = "D"
$ws3.Range("E1").Value = "E"
$ws3.Range("F1").Value = "F"
